# [PBL6] Range F0 from 60 to 400 Hz for cepstral method
# Rewrites the F0 comparison data (WaveSurfer ref vs cepstrum-based F0) for the new
# 60-400 Hz analysis range: one fewer data row (61 -> 60 rows), and updated values
# for columns A (F0_WaveSurfer), B (F0_Cepstrum based pitch), C (Result = |A-B|) and
# the D2 AVERAGE summary cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The recomputed dataset has 59 data rows (rows 2-60) instead of 60 (rows 2-61),
# so drop the now-unused last row first.
$ws.Rows.Item(61).Delete()

$ws.Cells.Item(2, 1).Value = 130.5646514892578
$ws.Cells.Item(2, 2).Value = 116.79
$ws.Cells.Item(2, 3).Value = 13.77465148925781

$ws.Cells.Item(3, 1).Value = 117.630859375
$ws.Cells.Item(3, 2).Value = 116.79
$ws.Cells.Item(3, 3).Value = 0.8408593749999937

$ws.Cells.Item(4, 1).Value = 117.5756530761719
$ws.Cells.Item(4, 2).Value = 117.65
$ws.Cells.Item(4, 3).Value = 0.07434692382813068

$ws.Cells.Item(5, 1).Value = 116.8266830444336
$ws.Cells.Item(5, 2).Value = 117.65
$ws.Cells.Item(5, 3).Value = 0.8233169555664119

$ws.Cells.Item(6, 1).Value = 117.0158767700195
$ws.Cells.Item(6, 2).Value = 115.94
$ws.Cells.Item(6, 3).Value = 1.075876770019534

$ws.Cells.Item(7, 1).Value = 118.3009262084961
$ws.Cells.Item(7, 2).Value = 117.65
$ws.Cells.Item(7, 3).Value = 0.6509262084960881

$ws.Cells.Item(8, 1).Value = 119.2888031005859
$ws.Cells.Item(8, 2).Value = 119.4
$ws.Cells.Item(8, 3).Value = 0.1111968994140682

$ws.Cells.Item(9, 1).Value = 120.0970840454102
$ws.Cells.Item(9, 2).Value = 119.4
$ws.Cells.Item(9, 3).Value = 0.6970840454101506

$ws.Cells.Item(10, 1).Value = 119.5390014648438
$ws.Cells.Item(10, 2).Value = 120.3
$ws.Cells.Item(10, 3).Value = 0.7609985351562472

$ws.Cells.Item(11, 1).Value = 117.8120651245117
$ws.Cells.Item(11, 2).Value = 117.65
$ws.Cells.Item(11, 3).Value = 0.1620651245117131

$ws.Cells.Item(12, 1).Value = 116.4637832641602
$ws.Cells.Item(12, 2).Value = 116.79
$ws.Cells.Item(12, 3).Value = 0.32621673583985

$ws.Cells.Item(13, 1).Value = 115.3791580200195
$ws.Cells.Item(13, 2).Value = 115.11
$ws.Cells.Item(13, 3).Value = 0.2691580200195318

$ws.Cells.Item(14, 1).Value = 113.8033447265625
$ws.Cells.Item(14, 2).Value = 114.29
$ws.Cells.Item(14, 3).Value = 0.4866552734375063

$ws.Cells.Item(15, 1).Value = 113.0747375488281
$ws.Cells.Item(15, 2).Value = 113.48
$ws.Cells.Item(15, 3).Value = 0.405262451171879

$ws.Cells.Item(16, 1).Value = 0
$ws.Cells.Item(16, 2).Value = 253.97
$ws.Cells.Item(16, 3).Value = 253.97

$ws.Cells.Item(17, 1).Value = 114.9552001953125
$ws.Cells.Item(17, 2).Value = 115.11
$ws.Cells.Item(17, 3).Value = 0.1547998046874994

$ws.Cells.Item(18, 1).Value = 114.5822982788086
$ws.Cells.Item(18, 2).Value = 116.79
$ws.Cells.Item(18, 3).Value = 2.207701721191413

$ws.Cells.Item(19, 1).Value = 113.5789642333984
$ws.Cells.Item(19, 2).Value = 113.48
$ws.Cells.Item(19, 3).Value = 0.09896423339843352

$ws.Cells.Item(20, 1).Value = 112.8562088012695
$ws.Cells.Item(20, 2).Value = 113.48
$ws.Cells.Item(20, 3).Value = 0.6237911987304727

$ws.Cells.Item(21, 1).Value = 113.127326965332
$ws.Cells.Item(21, 2).Value = 113.48
$ws.Cells.Item(21, 3).Value = 0.3526730346679727

$ws.Cells.Item(22, 1).Value = 114.2129516601562
$ws.Cells.Item(22, 2).Value = 114.29
$ws.Cells.Item(22, 3).Value = 0.07704833984375625

$ws.Cells.Item(23, 1).Value = 114.3563995361328
$ws.Cells.Item(23, 2).Value = 114.29
$ws.Cells.Item(23, 3).Value = 0.06639953613280625

$ws.Cells.Item(24, 1).Value = 114.4540557861328
$ws.Cells.Item(24, 2).Value = 115.11
$ws.Cells.Item(24, 3).Value = 0.6559442138671869

$ws.Cells.Item(25, 1).Value = 116.7946014404297
$ws.Cells.Item(25, 2).Value = 115.94
$ws.Cells.Item(25, 3).Value = 0.8546014404296898

$ws.Cells.Item(26, 1).Value = 117.8669128417969
$ws.Cells.Item(26, 2).Value = 117.65
$ws.Cells.Item(26, 3).Value = 0.2169128417968693

$ws.Cells.Item(27, 1).Value = 118.1022720336914
$ws.Cells.Item(27, 2).Value = 118.52
$ws.Cells.Item(27, 3).Value = 0.4177279663085898

$ws.Cells.Item(28, 1).Value = 117.1706085205078
$ws.Cells.Item(28, 2).Value = 117.65
$ws.Cells.Item(28, 3).Value = 0.4793914794921932

$ws.Cells.Item(29, 1).Value = 115.7100830078125
$ws.Cells.Item(29, 2).Value = 115.94
$ws.Cells.Item(29, 3).Value = 0.2299169921874977

$ws.Cells.Item(30, 1).Value = 114.5855102539062
$ws.Cells.Item(30, 2).Value = 114.29
$ws.Cells.Item(30, 3).Value = 0.2955102539062437

$ws.Cells.Item(31, 1).Value = 0
$ws.Cells.Item(31, 2).Value = 114.29
$ws.Cells.Item(31, 3).Value = 114.29

$ws.Cells.Item(32, 1).Value = 124.789680480957
$ws.Cells.Item(32, 2).Value = 124.03
$ws.Cells.Item(32, 3).Value = 0.7596804809570301

$ws.Cells.Item(33, 1).Value = 122.5589828491211
$ws.Cells.Item(33, 2).Value = 121.21
$ws.Cells.Item(33, 3).Value = 1.3489828491211

$ws.Cells.Item(34, 1).Value = 123.3526458740234
$ws.Cells.Item(34, 2).Value = 123.08
$ws.Cells.Item(34, 3).Value = 0.2726458740234392

$ws.Cells.Item(35, 1).Value = 123.6247711181641
$ws.Cells.Item(35, 2).Value = 124.03
$ws.Cells.Item(35, 3).Value = 0.4052288818359386

$ws.Cells.Item(36, 1).Value = 123.3014678955078
$ws.Cells.Item(36, 2).Value = 123.08
$ws.Cells.Item(36, 3).Value = 0.2214678955078142

$ws.Cells.Item(37, 1).Value = 122.4353866577148
$ws.Cells.Item(37, 2).Value = 122.14
$ws.Cells.Item(37, 3).Value = 0.2953866577148432

$ws.Cells.Item(38, 1).Value = 122.0922622680664
$ws.Cells.Item(38, 2).Value = 122.14
$ws.Cells.Item(38, 3).Value = 0.04773773193359432

$ws.Cells.Item(39, 1).Value = 119.9362945556641
$ws.Cells.Item(39, 2).Value = 120.3
$ws.Cells.Item(39, 3).Value = 0.3637054443359347

$ws.Cells.Item(40, 1).Value = 119.4052810668945
$ws.Cells.Item(40, 2).Value = 120.3
$ws.Cells.Item(40, 3).Value = 0.8947189331054659

$ws.Cells.Item(41, 1).Value = 118.4627532958984
$ws.Cells.Item(41, 2).Value = 118.52
$ws.Cells.Item(41, 3).Value = 0.05724670410155852

$ws.Cells.Item(42, 1).Value = 118.5424423217773
$ws.Cells.Item(42, 2).Value = 120.3
$ws.Cells.Item(42, 3).Value = 1.757557678222653

$ws.Cells.Item(43, 1).Value = 115.0749282836914
$ws.Cells.Item(43, 2).Value = 115.94
$ws.Cells.Item(43, 3).Value = 0.8650717163085915

$ws.Cells.Item(44, 1).Value = 115.3555145263672
$ws.Cells.Item(44, 2).Value = 115.11
$ws.Cells.Item(44, 3).Value = 0.2455145263671881

$ws.Cells.Item(45, 1).Value = 114.7453155517578
$ws.Cells.Item(45, 2).Value = 115.11
$ws.Cells.Item(45, 3).Value = 0.3646844482421869

$ws.Cells.Item(46, 1).Value = 113.8243560791016
$ws.Cells.Item(46, 2).Value = 113.48
$ws.Cells.Item(46, 3).Value = 0.3443560791015585

$ws.Cells.Item(47, 1).Value = 0
$ws.Cells.Item(47, 2).Value = 228.57
$ws.Cells.Item(47, 3).Value = 228.57

$ws.Cells.Item(48, 1).Value = 112.0472793579102
$ws.Cells.Item(48, 2).Value = 200
$ws.Cells.Item(48, 3).Value = 87.95272064208984

$ws.Cells.Item(49, 1).Value = 114.4084854125977
$ws.Cells.Item(49, 2).Value = 115.11
$ws.Cells.Item(49, 3).Value = 0.7015145874023432

$ws.Cells.Item(50, 1).Value = 114.8593978881836
$ws.Cells.Item(50, 2).Value = 115.11
$ws.Cells.Item(50, 3).Value = 0.2506021118164057

$ws.Cells.Item(51, 1).Value = 113.8720703125
$ws.Cells.Item(51, 2).Value = 115.11
$ws.Cells.Item(51, 3).Value = 1.237929687499999

$ws.Cells.Item(52, 1).Value = 0
$ws.Cells.Item(52, 2).Value = 124.03
$ws.Cells.Item(52, 3).Value = 124.03

$ws.Cells.Item(53, 1).Value = 112.0294876098633
$ws.Cells.Item(53, 2).Value = 111.11
$ws.Cells.Item(53, 3).Value = 0.9194876098632818

$ws.Cells.Item(54, 1).Value = 114.5335006713867
$ws.Cells.Item(54, 2).Value = 112.68
$ws.Cells.Item(54, 3).Value = 1.853500671386712

$ws.Cells.Item(55, 1).Value = 118.2854843139648
$ws.Cells.Item(55, 2).Value = 116.79
$ws.Cells.Item(55, 3).Value = 1.495484313964837

$ws.Cells.Item(56, 1).Value = 128.6198425292969
$ws.Cells.Item(56, 2).Value = 130.08
$ws.Cells.Item(56, 3).Value = 1.460157470703138

$ws.Cells.Item(57, 1).Value = 116.5141220092773
$ws.Cells.Item(57, 2).Value = 115.94
$ws.Cells.Item(57, 3).Value = 0.574122009277346

$ws.Cells.Item(58, 1).Value = 115.5278549194336
$ws.Cells.Item(58, 2).Value = 116.79
$ws.Cells.Item(58, 3).Value = 1.262145080566413

$ws.Cells.Item(59, 1).Value = 112.9555969238281
$ws.Cells.Item(59, 2).Value = 111.11
$ws.Cells.Item(59, 3).Value = 1.845596923828126

$ws.Cells.Item(60, 1).Value = 0
$ws.Cells.Item(60, 2).Value = 94.12
$ws.Cells.Item(60, 3).Value = 94.12

# Summary average of the Result column over the new row range (C2:C60)
$ws.Cells.Item(2, 4).Value = 16.10113974361096
